# Scheduled-runner style refresh of market/profit figures across the
# per-job "Leve" sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR). Only the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) move;
# columns A:G (name/item/level/exp/gil/amount/itemId) are untouched.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values,
        [string[]]$ClearCols = @()
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Cells.Item($Row, $col).Value = $Values[$col]
    }
    foreach ($col in $ClearCols) {
        $ws.Cells.Item($Row, $col).ClearContents()
    }
}

# Column letters -> 1-based column index (A=1 ... N=14)
# H=8 I=9 J=10 K=11 L=12 M=13 N=14

# ---------------- ALC ----------------
Set-Row -SheetName "ALC" -Row 4   -Values @{8=2509.6; 9=1683.3334; 11=1683.3334; 13=-1569.3334}
Set-Row -SheetName "ALC" -Row 15  -Values @{8=731.9459000000001; 9=731.9459000000001; 11=2195.8377; 13=-2026.8377}
Set-Row -SheetName "ALC" -Row 99  -Values @{8=526.5; 10=839.8; 12=2519.4; 14=-5515.4}
Set-Row -SheetName "ALC" -Row 111 -Values @{8=2488.0454; 9=2176.2666; 10=3156.1428; 11=6528.7998; 12=9468.428400000001; 13=-3461.7998; 14=-15602.4284}
Set-Row -SheetName "ALC" -Row 112 -Values @{8=1891.7142; 9=1324.5; 10=2118.6; 11=3973.5; 12=6355.799999999999; 13=-2865.5; 14=-8571.799999999999}
Set-Row -SheetName "ALC" -Row 116 -Values @{8=12756.529; 9=17598.8; 10=11921.655; 11=17598.8; 12=11921.655; 13=-14156.8; 14=-18805.655}
Set-Row -SheetName "ALC" -Row 141 -Values @{8=7040.4443; 9=3604.3635; 11=10813.0905; 13=-5633.0905}

# ---------------- ARM ----------------
Set-Row -SheetName "ARM" -Row 2   -Values @{8=3475.2; 9=3594; 10=3000; 11=3594; 12=3000; 13=-3481; 14=-3226}
Set-Row -SheetName "ARM" -Row 36  -Values @{8=3150; 9=3150; 11=3150; 13=-2804}
Set-Row -SheetName "ARM" -Row 45  -Values @{8=3496.1; 9=1396; 10=5596.2; 11=1396; 12=5596.2; 13=-1019; 14=-6350.2}
Set-Row -SheetName "ARM" -Row 102 -Values @{8=1886.3529; 9=1886.3529; 11=1886.3529; 13=-264.3529000000001}
Set-Row -SheetName "ARM" -Row 116 -Values @{8=3475.2; 9=3594; 10=3000; 11=3594; 12=3000; 13=-1300; 14=-7588}

# ---------------- BSM ----------------
Set-Row -SheetName "BSM" -Row 3   -Values @{8=3475.2; 9=3594; 10=3000; 11=3594; 12=3000; 13=-3480; 14=-3228}
Set-Row -SheetName "BSM" -Row 20  -Values @{8=1174.2122; 9=935.5714; 11=935.5714; 13=-688.5714}
Set-Row -SheetName "BSM" -Row 64  -Values @{8=4736.3; 9=10006; 10=3418.875; 11=10006; 12=3418.875; 13=-9781; 14=-3868.875}
Set-Row -SheetName "BSM" -Row 67  -Values @{8=4736.3; 9=10006; 10=3418.875; 11=10006; 12=3418.875; 13=-9226; 14=-4978.875}
Set-Row -SheetName "BSM" -Row 88  -Values @{8=22374.75; 10=22374.75; 12=22374.75; 14=-23186.75}
Set-Row -SheetName "BSM" -Row 91  -Values @{8=22374.75; 10=22374.75; 12=22374.75; 14=-25182.75}
Set-Row -SheetName "BSM" -Row 105 -Values @{8=2070.889; 9=1815.4166; 10=2581.8333; 11=1815.4166; 12=2581.8333; 13=-68.41660000000002; 14=-6075.8333}
Set-Row -SheetName "BSM" -Row 107 -Values @{8=1695.711; 10=1610.8334; 12=1610.8334; 14=-5450.8334}

# ---------------- CRP ----------------
Set-Row -SheetName "CRP" -Row 22  -Values @{8=543.25; 10=588.7143; 12=588.7143; 14=-1288.7143}
Set-Row -SheetName "CRP" -Row 23  -Values @{8=17500; 9=17500; 11=17500; 13=-17260}
Set-Row -SheetName "CRP" -Row 27  -Values @{8=17500; 9=17500; 11=17500; 13=-17308}
Set-Row -SheetName "CRP" -Row 31  -Values @{8=3257.1904; 9=1700.8334; 10=5332.3335; 11=1700.8334; 12=5332.3335; 13=-1405.8334; 14=-5922.3335}
Set-Row -SheetName "CRP" -Row 34  -Values @{8=3257.1904; 9=1700.8334; 10=5332.3335; 11=1700.8334; 12=5332.3335; 13=-1498.8334; 14=-5736.3335}
Set-Row -SheetName "CRP" -Row 105 -Values @{8=1937.2084; 9=1621.875; 11=1621.875; 13=125.125}
Set-Row -SheetName "CRP" -Row 107 -Values @{8=1595.8823; 9=1596.5; 11=1596.5; 13=323.5}
Set-Row -SheetName "CRP" -Row 134 -Values @{8=4225.4375; 10=8502.666999999999; 12=25508.001; 14=-30578.001}
Set-Row -SheetName "CRP" -Row 141 -Values @{8=76447; 10=76447; 12=76447; 14=-86807}

# ---------------- GSM ----------------
# Row 53 also drops the LeveProfitNQ (M) cell entirely - the NQ price now
# matches the sale price exactly so profit/loss is only tracked on the N column.
Set-Row -SheetName "GSM" -Row 53  -Values @{8=49999; 9=0; 10=49999; 11=0; 12=49999; 14=-51261} -ClearCols @(13)
Set-Row -SheetName "GSM" -Row 113 -Values @{8=2155.7368; 9=1813.6154; 11=1813.6154; 13=356.3846000000001}

# ---------------- LTW ----------------
Set-Row -SheetName "LTW" -Row 61  -Values @{8=2329.6667; 9=2329.6667; 11=2329.6667; 13=-2127.6667}
Set-Row -SheetName "LTW" -Row 113 -Values @{8=2329.6667; 9=2329.6667; 11=2329.6667; 13=-159.6667000000002}

# ---------------- WVR ----------------
Set-Row -SheetName "WVR" -Row 96  -Values @{8=3057.2; 9=1900; 10=3553.1428; 11=1900; 12=3553.1428; 13=-527; 14=-6299.1428}
Set-Row -SheetName "WVR" -Row 113 -Values @{8=356.57144; 9=349.33334; 11=1048.00002; 13=1121.99998}

Write-Output "Done."
